# Rename the inline picture shapes that live in the headers/footers.
#
#   footer "first"   (id=3) : image2.png -> image1.png   (Pearson logo)
#   footer "default"  (id=2) : image2.png -> image1.png   (Pearson logo)
#   header "first"    (id=1) : image1.jpg -> image2.jpg   (BTEC logo)
#
# Both the Pearson logo inline shapes carry the descr
# "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png"
# and the BTEC logo inline shape carries the descr "BTec_Logo-Orange" -- use
# that to identify each picture unambiguously before renaming it.

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {

    # --- Headers (primary / first-page / even-page) ---
    for ($i = 1; $i -le 3; $i++) {
        $hdr = $sec.Headers.Item($i)
        if ($hdr.Exists) {
            foreach ($shp in $hdr.Range.InlineShapes) {
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                }
            }
        }
    }

    # --- Footers (primary / first-page / even-page) ---
    for ($i = 1; $i -le 3; $i++) {
        $ftr = $sec.Footers.Item($i)
        if ($ftr.Exists) {
            foreach ($shp in $ftr.Range.InlineShapes) {
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }
}
